$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.882.11'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.630.80'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.27'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.32'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.258'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.44%  '
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0880'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.862.17'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '1.628.62'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.41'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').Value = '27.878.18'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.54'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('E22').Value = '  -3.94%  '
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('E24').Value = '  -3.69%  '
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.91'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.62'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.995'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0481'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.09'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '1.400.71'
$ws.Range('E34').Value = '  -2.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.57'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.65%  '
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('E38').Value = '  +1.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.557'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.868'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.994'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '66.77'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.75%  '
$ws.Range('E44').Value = '  +1.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.45'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('E46').Value = '  -0.75%  '
$ws.Range('D47').Value = '1.772.20'
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.11'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('E49').Value = '  -4.07%  '
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('E51').Value = '  -0.11%  '
